$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume snapshot (GitHub Actions refresh).
# For cells whose new text looks like a plain number (e.g. "1.00", "0.0514"),
# force the cell to Text format first so Excel keeps the exact original
# string (preserving trailing zeros / leading zeros) instead of silently
# converting it to a numeric value; the format is reset back to Normal
# immediately afterwards so no stray formatting is left behind.

$ws.Range("D2").Value = "34.705.72"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "1.788.96"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.43%  "
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("E10").Value = "  +2.93%  "
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").Value = "2.047.19"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.63%  "
$ws.Range("D14").Value = "1.781.12"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "34.707.89"
$ws.Range("E15").Value = "  +2.67%  "
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("E17").Value = "  +3.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("E20").Value = "  +5.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0514"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "1.440.57"
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.633"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0190"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "82.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0503"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("D47").Value = "1.942.21"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  -2.48%  "
